$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.57800190790492
$ws.Range("C2").Value = 8.000055910654073
$ws.Range("D2").Value = 7.940792452306904
$ws.Range("E2").Value = 12.99624494455836
$ws.Range("F2").Value = 38.13641248422841
$ws.Range("J2").Value = 10.24051627375574
$ws.Range("K2").Value = 11.68610621854213
$ws.Range("L2").Value = 10.7926058205909
$ws.Range("N2").Value = 21.20514023523401
$ws.Range("O2").Value = 29.45693456936716
$ws.Range("B3").Value = 15.40309458605452
$ws.Range("C3").Value = 7.979124672057908
$ws.Range("D3").Value = 7.925965987496625
$ws.Range("E3").Value = 13.01381980262542
$ws.Range("F3").Value = 38.21399427253787
$ws.Range("J3").Value = 10.26006266581307
$ws.Range("K3").Value = 11.55980489088751
$ws.Range("L3").Value = 10.79251960207061
$ws.Range("N3").Value = 21.26589487856848
$ws.Range("O3").Value = 29.53561375106056
$ws.Range("B4").Value = 15.29757681231026
$ws.Range("C4").Value = 7.966175663296172
$ws.Range("D4").Value = 7.917934746494919
$ws.Range("E4").Value = 13.02618630448242
$ws.Range("F4").Value = 38.26871890067844
$ws.Range("J4").Value = 10.2728155832853
$ws.Range("K4").Value = 11.48348350119378
$ws.Range("L4").Value = 10.79376169256781
$ws.Range("N4").Value = 21.30495113743897
$ws.Range("O4").Value = 29.58887585530327
$ws.Range("B5").Value = 15.25509596058923
$ws.Range("C5").Value = 7.960875547451754
$ws.Range("D5").Value = 7.914934082723703
$ws.Range("E5").Value = 13.03162245094352
$ws.Range("F5").Value = 38.29280074596404
$ws.Range("J5").Value = 10.27820188468306
$ws.Range("K5").Value = 11.45272296537458
$ws.Range("L5").Value = 10.79459463752374
$ws.Range("N5").Value = 21.32130888670471
$ws.Range("O5").Value = 29.61182520156541
$ws.Range("B6").Value = 15.24807471192447
$ws.Range("C6").Value = 7.959994093636411
$ws.Range("D6").Value = 7.914452330602836
$ws.Range("E6").Value = 13.03254909328347
$ws.Range("F6").Value = 38.29690704390661
$ws.Range("J6").Value = 10.27910772751406
$ws.Range("K6").Value = 11.44763673052567
$ws.Range("L6").Value = 10.79475271581459
$ws.Range("N6").Value = 21.32405181245001
$ws.Range("O6").Value = 29.61571106933761
$ws.Range("B7").Value = 15.297001740215
$ws.Range("C7").Value = 7.966104276882775
$ws.Range("D7").Value = 7.917893173458233
$ws.Range("E7").Value = 13.02625801136184
$ws.Range("F7").Value = 38.26903646738656
$ws.Range("J7").Value = 10.27288745742092
$ws.Range("K7").Value = 11.48306723081731
$ws.Range("L7").Value = 10.79377160138138
$ws.Range("N7").Value = 21.30516995239516
$ws.Range("O7").Value = 29.58918031997271
$ws.Range("B8").Value = 15.51733279963883
$ws.Range("C8").Value = 7.992859664414732
$ws.Range("D8").Value = 7.935459288182171
$ws.Range("E8").Value = 13.0019780964357
$ws.Range("F8").Value = 38.16169052082238
$ws.Range("J8").Value = 10.24710019703286
$ws.Range("K8").Value = 11.6423221195092
$ws.Range("L8").Value = 10.79230801517827
$ws.Range("N8").Value = 21.22572544864287
$ws.Range("O8").Value = 29.48303473815577
$ws.Range("B9").Value = 15.96219757783857
$ws.Range("C9").Value = 8.044507118706996
$ws.Range("D9").Value = 7.978297886676544
$ws.Range("E9").Value = 12.96684177193302
$ws.Range("F9").Value = 38.00749257953472
$ws.Range("J9").Value = 10.20247279943884
$ws.Range("K9").Value = 11.9629270720053
$ws.Range("L9").Value = 10.79966095780772
$ws.Range("N9").Value = 21.08378191923783
$ws.Range("O9").Value = 29.31421605821458
$ws.Range("B10").Value = 16.29401551694192
$ws.Range("C10").Value = 8.081888523913674
$ws.Range("D10").Value = 8.014721554801083
$ws.Range("E10").Value = 12.94859958686888
$ws.Range("F10").Value = 37.92859643713337
$ws.Range("J10").Value = 10.17327893401039
$ws.Range("K10").Value = 12.20159610485682
$ws.Range("L10").Value = 10.81122133576632
$ws.Range("N10").Value = 20.98785243877131
$ws.Range("O10").Value = 29.21420772065408
$ws.Range("B11").Value = 16.44544881833577
$ws.Range("C11").Value = 8.098758029753647
$ws.Range("D11").Value = 8.032327841286602
$ws.Range("E11").Value = 12.94193712707576
$ws.Range("F11").Value = 37.90018083472031
$ws.Range("J11").Value = 10.16077232783957
$ws.Range("K11").Value = 12.31043777752725
$ws.Range("L11").Value = 10.81780003650485
$ws.Range("N11").Value = 20.94600800697858
$ws.Range("O11").Value = 29.17393582521685
$ws.Range("B12").Value = 16.50281321677934
$ws.Range("C12").Value = 8.105125354075273
$ws.Range("D12").Value = 8.039140452357763
$ws.Range("E12").Value = 12.9396487363682
$ws.Range("F12").Value = 37.89049548342529
$ws.Range("J12").Value = 10.15614722683079
$ws.Range("K12").Value = 12.35165793605704
$ws.Range("L12").Value = 10.82047933528648
$ws.Range("N12").Value = 20.93041929946628
$ws.Range("O12").Value = 29.15943743722778
$ws.Range("B13").Value = 16.49045867286308
$ws.Range("C13").Value = 8.103754982502748
$ws.Range("D13").Value = 8.037666823172019
$ws.Range("E13").Value = 12.94013116301619
$ws.Range("F13").Value = 37.89253358631593
$ws.Range("J13").Value = 10.15713839996863
$ws.Range("K13").Value = 12.34278079694457
$ws.Range("L13").Value = 10.81989396323693
$ws.Range("N13").Value = 20.93376520047165
$ws.Range("O13").Value = 29.16252648160055
$ws.Range("B14").Value = 16.45016810784024
$ws.Range("C14").Value = 8.099282297580134
$ws.Range("D14").Value = 8.032885424310097
$ws.Range("E14").Value = 12.94174416357159
$ws.Range("F14").Value = 37.89936247029594
$ws.Range("J14").Value = 10.16038959806542
$ws.Range("K14").Value = 12.31382909740166
$ws.Range("L14").Value = 10.8180167045638
$ws.Range("N14").Value = 20.94472037315132
$ws.Range("O14").Value = 29.17272796498247
$ws.Range("B15").Value = 16.42549008805228
$ws.Range("C15").Value = 8.096539902252434
$ws.Range("D15").Value = 8.029975516712755
$ws.Range("E15").Value = 12.94276269471396
$ws.Range("F15").Value = 37.90368535742245
$ws.Range("J15").Value = 10.16239547884158
$ws.Range("K15").Value = 12.29609487130363
$ws.Range("L15").Value = 10.81689127181918
$ws.Range("N15").Value = 20.95146415071351
$ws.Range("O15").Value = 29.17907458232469
$ws.Range("B16").Value = 16.28412449476253
$ws.Range("C16").Value = 8.080783224074796
$ws.Range("D16").Value = 8.013591502973135
$ws.Range("E16").Value = 12.94906785461167
$ws.Range("F16").Value = 37.93060389575069
$ws.Range("J16").Value = 10.17411180758327
$ws.Range("K16").Value = 12.19448550423438
$ws.Range("L16").Value = 10.81081781721935
$ws.Range("N16").Value = 20.99062307891049
$ws.Range("O16").Value = 29.21694473631464
$ws.Range("B17").Value = 16.19749019103555
$ws.Range("C17").Value = 8.071081383499742
$ws.Range("D17").Value = 8.003803507023239
$ws.Range("E17").Value = 12.95335439002656
$ws.Range("F17").Value = 37.94903222165523
$ws.Range("J17").Value = 10.18149730871435
$ws.Range("K17").Value = 12.13219577694974
$ws.Range("L17").Value = 10.80742880871954
$ws.Range("N17").Value = 21.01510456715532
$ws.Range("O17").Value = 29.24151498007436
$ws.Range("B18").Value = 16.14770929310239
$ws.Range("C18").Value = 8.065488468057309
$ws.Range("D18").Value = 7.998271542692235
$ws.Range("E18").Value = 12.95597392587005
$ws.Range("F18").Value = 37.96033526843879
$ws.Range("J18").Value = 10.18581811357356
$ws.Range("K18").Value = 12.09639582462695
$ws.Range("L18").Value = 10.80560387699655
$ws.Range("N18").Value = 21.02935463477063
$ws.Range("O18").Value = 29.25613865344921
$ws.Range("B19").Value = 16.13086427232801
$ws.Range("C19").Value = 8.063592666060229
$ws.Range("D19").Value = 7.996415425240352
$ws.Range("E19").Value = 12.95688733226641
$ws.Range("F19").Value = 37.96428311201259
$ws.Range("J19").Value = 10.18729358977139
$ws.Range("K19").Value = 12.08428035467549
$ws.Range("L19").Value = 10.80500739162362
$ws.Range("N19").Value = 21.03420851784804
$ws.Range("O19").Value = 29.2611743730867
$ws.Range("B20").Value = 16.20670787001576
$ws.Range("C20").Value = 8.072115482091727
$ws.Range("D20").Value = 8.004835356417219
$ws.Range("E20").Value = 12.95288214411731
$ws.Range("F20").Value = 37.94699767655111
$ws.Range("J20").Value = 10.18070357166488
$ws.Range("K20").Value = 12.13882403626492
$ws.Range("L20").Value = 10.80777671814026
$ws.Range("N20").Value = 21.01248099230399
$ws.Range("O20").Value = 29.23884855970674
$ws.Range("B21").Value = 16.46200227953015
$ws.Range("C21").Value = 8.100596608312708
$ws.Range("D21").Value = 8.034285918068758
$ws.Range("E21").Value = 12.94126402685412
$ws.Range("F21").Value = 37.89732748616566
$ws.Range("J21").Value = 10.15943163628035
$ws.Range("K21").Value = 12.32233307316294
$ws.Range("L21").Value = 10.81856301034089
$ws.Range("N21").Value = 20.94149561115323
$ws.Range("O21").Value = 29.16971113562371
$ws.Range("B22").Value = 16.62894525478698
$ws.Range("C22").Value = 8.119088868944207
$ws.Range("D22").Value = 8.054379817317361
$ws.Range("E22").Value = 12.93503761807566
$ws.Range("F22").Value = 37.87113101478872
$ws.Range("J22").Value = 10.14617530531541
$ws.Range("K22").Value = 12.44227556295368
$ws.Range("L22").Value = 10.82670793672505
$ws.Range("N22").Value = 20.89659927037045
$ws.Range("O22").Value = 29.12890761412937
$ws.Range("B23").Value = 16.53985246662524
$ws.Range("C23").Value = 8.109230781076548
$ws.Range("D23").Value = 8.043579128869697
$ws.Range("E23").Value = 12.93823596520036
$ws.Range("F23").Value = 37.88453927226544
$ws.Range("J23").Value = 10.15319146957421
$ws.Range("K23").Value = 12.37827048962743
$ws.Range("L23").Value = 10.82226118864851
$ws.Range("N23").Value = 20.92042472306072
$ws.Range("O23").Value = 29.15028407159727
$ws.Range("B24").Value = 16.20254047152743
$ws.Range("C24").Value = 8.071648013365932
$ws.Range("D24").Value = 8.004368560242856
$ws.Range("E24").Value = 12.9530951631726
$ws.Range("F24").Value = 37.94791528792364
$ws.Range("J24").Value = 10.18106218721662
$ws.Range("K24").Value = 12.13582736255845
$ws.Range("L24").Value = 10.80761904362062
$ws.Range("N24").Value = 21.01366656424021
$ws.Range("O24").Value = 29.24005249743362
$ws.Range("B25").Value = 15.84077124276541
$ws.Range("C25").Value = 8.030629455865929
$ws.Range("D25").Value = 7.965826951853071
$ws.Range("E25").Value = 12.97501467847139
$ws.Range("F25").Value = 38.04317100009714
$ws.Range("J25").Value = 10.21391255347445
$ws.Range("K25").Value = 11.87550603525665
$ws.Range("L25").Value = 10.79658461531679
$ws.Range("N25").Value = 21.12070782624115
$ws.Range("O25").Value = 29.16252648160055

Write-Output "Updated 240 cells"
